# Update bank statement #198 to next statement period (statement_198.xlsx diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: account holder name / card number
$ws.Range("C2").Value = "Hartmut"
# Leading apostrophe forces the numeric-looking card number to stay text
# (matches the source cell, which is stored as text, not a number).
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 11.10.2023"

# Row 6
$ws.Range("B6").Value = "13.10."
$ws.Range("C6").Value = "14.10."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-31591718"
$ws.Range("E6").Value = "56,44-"

# Row 7
$ws.Range("B7").Value = "14.10."
$ws.Range("C7").Value = "15.10."
$ws.Range("D7").Value = "MCDONALDS Rastatt"
$ws.Range("E7").Value = "44,29-"

# Row 8
$ws.Range("B8").Value = "16.10."
$ws.Range("C8").Value = "17.10."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 82423078"
$ws.Range("E8").Value = "39,25-"

# Row 9 (newly populated; align style with rows 6-8/12 for the amount column)
$ws.Range("B9").Value = "19.10."
$ws.Range("C9").Value = "20.10."
$ws.Range("D9").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E9").Value = "78,89-"
$ws.Range("E9").HorizontalAlignment = -4152
$ws.Range("E9").VerticalAlignment = -4107
$ws.Range("E9").WrapText = $false

# Row 10 (newly populated)
$ws.Range("B10").Value = "23.10."
$ws.Range("C10").Value = "24.10."
$ws.Range("D10").Value = "AMAZON.DE MKTPLC EU RYBIHN"
$ws.Range("E10").Value = "59,74-"
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4107
$ws.Range("E10").WrapText = $false

# Row 11 (newly populated)
$ws.Range("B11").Value = "26.10."
$ws.Range("C11").Value = "27.10."
$ws.Range("D11").Value = "BURGER KING Ludwigsburg"
$ws.Range("E11").Value = "11,62-"
$ws.Range("E11").HorizontalAlignment = -4152
$ws.Range("E11").VerticalAlignment = -4107
$ws.Range("E11").WrapText = $false

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 28.10.2023"
$ws.Range("E12").Value = "290,23-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 02.11.2023"
